$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Email" header in G1 and the email value in G2
$ws.Range("G1").Value = "Email"
$ws.Range("G2").Value = "adrianrentea01@gmail.com"

# Style the email cell like a hyperlink: underlined, blue font
$ws.Range("G2").Font.Underline = 2
$ws.Range("G2").Font.Color = 13395456
$ws.Range("G2").Font.Name = "Calibri"
$ws.Range("G2").Font.Size = 12

# Update selection to match the new target cell(s)
$ws.Range("G1:G2").Select()
